$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update data values
#    - cap_factors!C2: 0.2574799786438868 -> 0.7
#    - demand!B2:      463.9011865031144  -> 500
# ---------------------------------------------------------------
$wsCapFactors = $wb.Worksheets.Item("cap_factors")
$wsCapFactors.Range("C2").Value = 0.7

$wsDemand = $wb.Worksheets.Item("demand")
$wsDemand.Range("B2").Value = 500

# ---------------------------------------------------------------
# 2) Highlight the input cells on config / thermal / vres with the
#    same accent fill colour (a light green "Accent 6" shade).
# ---------------------------------------------------------------
$wsConfig = $wb.Worksheets.Item("config")
$wsConfig.Range("B1:B2").Interior.ThemeColor = 10

$wsThermal = $wb.Worksheets.Item("thermal")
$wsThermal.Range("B2:F2").Interior.ThemeColor = 10

$wsVres = $wb.Worksheets.Item("vres")
$wsVres.Range("B2:D2").Interior.ThemeColor = 10

# ---------------------------------------------------------------
# 3) Hide the ramp-rate columns (D:E) on the thermal sheet.
# ---------------------------------------------------------------
$wsThermal.Range("D1:E1").EntireColumn.Hidden = $true

# ---------------------------------------------------------------
# 4) Reorder the tabs so "cap_factors" comes before "demand".
# ---------------------------------------------------------------
$wsCapFactors.Move($wsDemand)

# ---------------------------------------------------------------
# 5) Update the selections shown on each sheet.
# ---------------------------------------------------------------
$wsConfig.Range("B1:B2").Select()
$wsThermal.Range("F2").Select()
$wsVres.Range("E44").Select()

# ---------------------------------------------------------------
# 6) Make "demand" the active tab/sheet (it is now the last tab).
#    Re-fetch it by name since its tab position changed in step 4.
# ---------------------------------------------------------------
$wsDemand = $wb.Worksheets.Item("demand")
$wsDemand.Activate()
